# enemyDatabase.xlsx edit:
#  - Fixed spellbook arrow layer (ATK/weakness-resist values on the
#    existing Doppelganger BLUE/YELLOW/RED rows, and max_hp bump to 150).
#  - Added Doppelganger AI: a brand-new "Doppelganger (???)" enemy row,
#    and re-tuned the existing "Doppelganger (GRAY)" row, which is pushed
#    one row down to make room for the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Doppelganger (BLUE) -- max_hp 100 -> 150, bolt wk/res 1 -> 0.5
$ws.Cells.Item(12, 4).Value = 150
$ws.Cells.Item(12, 17).Value = 0.5

# --- Row 13: Doppelganger (YELLOW) -- max_hp 100 -> 150,
#     agni/cryo/bolt wk-res shifted from (2, -1, 1) to (0.5, 2, -1)
$ws.Cells.Item(13, 4).Value = 150
$ws.Cells.Item(13, 15).Value = 0.5
$ws.Cells.Item(13, 16).Value = 2
$ws.Cells.Item(13, 17).Value = -1

# --- Row 14: Doppelganger (RED) -- max_hp 100 -> 150,
#     agni/cryo/bolt wk-res shifted from (2, -1, 1) to (-1, 0.5, 2)
$ws.Cells.Item(14, 4).Value = 150
$ws.Cells.Item(14, 15).Value = -1
$ws.Cells.Item(14, 16).Value = 0.5
$ws.Cells.Item(14, 17).Value = 2

# --- Insert a brand-new row for "Doppelganger (???)" above the
#     existing "Doppelganger (GRAY)" row (pushes GRAY + END down by one).
$ws.Rows.Item(15).Insert()

# New row 15: Doppelganger (???)
$ws.Cells.Item(15, 1).Value = "Doppelganger (???)"
$ws.Cells.Item(15, 2).Value = "spr_bt_doppelganger_b_placeholder"
$ws.Cells.Item(15, 3).Value = "doppelganger_1"
$ws.Cells.Item(15, 4).Value = 150
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(15, 7).Value = 1.5
$ws.Cells.Item(15, 8).Value = 0.25
$ws.Cells.Item(15, 9).Value = 0.95
$ws.Cells.Item(15, 10).Value = 1
$ws.Cells.Item(15, 11).Value = 5
$ws.Cells.Item(15, 12).Value = "Doppleganger1"
$ws.Cells.Item(15, 13).Value = "none"
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
$ws.Cells.Item(15, 17).Value = 0
$ws.Cells.Item(15, 18).Value = "GROUP/DEFAULT"
$ws.Cells.Item(15, 19).Value = "sword"
$ws.Cells.Item(15, 20).Value = "null"
$ws.Cells.Item(15, 21).Value = "null"
$ws.Cells.Item(15, 22).Value = "END"

# Row 16 (previously row 15): Doppelganger (GRAY), re-tuned stats.
$ws.Cells.Item(16, 4).Value = 75
$ws.Cells.Item(16, 7).Value = 2
$ws.Cells.Item(16, 9).Value = 3
$ws.Cells.Item(16, 14).Value = 2
$ws.Cells.Item(16, 15).Value = 0
$ws.Cells.Item(16, 16).Value = 0
$ws.Cells.Item(16, 17).Value = 0
